$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rebuild the email hyperlink on the new "Email Address" column (E) ---
# Remove every existing hyperlink (old D2 -> mailto:test@test.com and old I2 -> image url).
$ws.Hyperlinks.Delete()

# --- Header row (row 1): rename / reposition columns ---
# A1 First Name, B1 Last Name, C1 Designation stay the same.
$ws.Range("D1").Value = "Company"
$ws.Range("E1").Value = "Email Address"
$ws.Range("F1").Value = "Speaker Bio"
$ws.Range("G1").Value = "Phone No"
$ws.Range("H1").Value = "Rating"

# --- Data row (row 2) ---
$ws.Range("D2").Value = "test"
$ws.Range("E2").Value = "test@test.com"
$ws.Range("F2").Value = "test"
$ws.Range("G2").Value = 12345667890
$ws.Range("H2").Value = "active / deactive"

# D2 lost its hyperlink (moved to E2), so drop its old blue/underline hyperlink
# formatting back to the plain/default look used by the rest of the row.
$ws.Range("A2").Copy()
$ws.Range("D2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Re-add the hyperlink for the email cell (now column E).
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:test@test.com", "", "", "test@test.com")

# Copy formatting (the existing hyperlink style) from the old "Profile Picture" link
# cell onto the new email cell so it keeps the blue hyperlink look / same style index,
# then remove the now-unused "Profile Picture" column entirely.
$ws.Range("I2").Copy()
$ws.Range("E2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("I1:I2").Clear()

# --- Selection ---
$null = $ws.Range("H2").Select()
